$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47:I47").Copy()
$ws.Range("A48:I48").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A48").Value = 46003
$ws.Range("B48").Value = 5610
$ws.Range("C48").Value = 4036
$ws.Range("D48").Value = 3741
$ws.Range("E48").Value = 219
$ws.Range("F48").Value = 47
$ws.Range("G48").Value = 28
$ws.Range("H48").Value = 1
$ws.Range("I48").Value = 0

$ws.Range("A48:I48").Select()
